# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计", with the
#    same per-fund layout used by the other quarterly sheets.
# 2. Insert a new leading row into "总计" summarising the 2022-Q1 quarter.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 1. New "2022-Q1" sheet, positioned right before "总计".
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

# Adding a sheet ahead of "总计" re-seats any previously captured handle to
# it (it silently starts pointing at the freshly inserted sheet instead) -
# grab a fresh handle by name now that the insert has happened.
$total = $wb.Worksheets.Item("总计")

# Pull over the existing quarterly layout/formatting (borders, bold header,
# centered alignment) so the new sheet's styles line up with its siblings
# instead of minting brand-new style entries. Column A's header cell (A1)
# is never populated on the sibling sheets either, so copy it separately
# from B1:H3 to avoid materialising a spurious empty A1 cell.
$q4.Range("B1:H3").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)
$q4.Range("A2:A3").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "288001"
$newSheet.Range("C2").Value = "华夏经典配置混合"
$newSheet.Range("D2").Value = "18.49"
$newSheet.Range("E2").Value = "63.85"
$newSheet.Range("F2").Value = "2.68"
$newSheet.Range("G2").Value = "0.4955"
$newSheet.Range("H2").Value = 10

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3:G3").NumberFormat = "@"
$newSheet.Range("B3").Value = "005347"
$newSheet.Range("C3").Value = "诺德量化优选6个月持有期混合"
$newSheet.Range("D3").Value = "2.60"
$newSheet.Range("E3").Value = "93.66"
$newSheet.Range("F3").Value = "2.76"
$newSheet.Range("G3").Value = "0.0718"
$newSheet.Range("H3").Value = 10

# The NumberFormat="@" nudge above was only needed so Excel wouldn't mangle
# the fund codes / decimal strings into numbers (dropping the leading zero
# on "005347"). Re-stamp the plain General formatting from the sibling
# sheet on top so the stored style matches (values stay text either way).
$q4.Range("B2:G3").Copy()
$newSheet.Range("B2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. New leading row in "总计" for the 2022-Q1 summary.
# ---------------------------------------------------------------------------
$total.Rows(2).Insert()
$total.Range("A2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.57

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

Write-Host "2022-Q1 sheet + 总计 summary row added"
